$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite rows 2-31 with the final data set (10 new rows inserted at top of data,
# original 20 data rows shifted down by 10, timestamps renumbered by 100 per row).

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "struggle"
$ws.Cells.Item(2, 3).Value = -0.1118526458740234
$ws.Cells.Item(2, 4).Value = 0.0269185900688171
$ws.Cells.Item(2, 5).Value = 0.0618541836738586
$ws.Cells.Item(2, 6).Value = -0.0058032199740409
$ws.Cells.Item(2, 7).Value = 0.0117591563612222
$ws.Cells.Item(2, 8).Value = 0.0122173046693205

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "struggle"
$ws.Cells.Item(3, 3).Value = -0.188694953918457
$ws.Cells.Item(3, 4).Value = -0.0127399563789367
$ws.Cells.Item(3, 5).Value = 0.0153613984584808
$ws.Cells.Item(3, 6).Value = -0.0114537235349416
$ws.Cells.Item(3, 7).Value = -0.0096211275085806
$ws.Cells.Item(3, 8).Value = -0.0482583530247211

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "struggle"
$ws.Cells.Item(4, 3).Value = -0.0261173248291015
$ws.Cells.Item(4, 4).Value = -0.1474769711494445
$ws.Cells.Item(4, 5).Value = 0.0655251443386077
$ws.Cells.Item(4, 6).Value = 0.0445931628346443
$ws.Cells.Item(4, 7).Value = 0.1122464910149574
$ws.Cells.Item(4, 8).Value = -0.0378736443817615

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "struggle"
$ws.Cells.Item(5, 3).Value = -0.1960973739624023
$ws.Cells.Item(5, 4).Value = 0.0549294650554657
$ws.Cells.Item(5, 5).Value = 0.0360765755176544
$ws.Cells.Item(5, 6).Value = 0.0612392425537109
$ws.Cells.Item(5, 7).Value = 0.09758572280406951
$ws.Cells.Item(5, 8).Value = -0.0021380283869802

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "struggle"
$ws.Cells.Item(6, 3).Value = -0.06610202789306641
$ws.Cells.Item(6, 4).Value = -0.1787786185741424
$ws.Cells.Item(6, 5).Value = 0.0745508223772049
$ws.Cells.Item(6, 6).Value = 0.0088575463742017
$ws.Cells.Item(6, 7).Value = 0.1237002089619636
$ws.Cells.Item(6, 8).Value = 0.0548251569271087

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "struggle"
$ws.Cells.Item(7, 3).Value = 0.0234136581420898
$ws.Cells.Item(7, 4).Value = 0.0270741879940032
$ws.Cells.Item(7, 5).Value = 0.2239813506603241
$ws.Cells.Item(7, 6).Value = -0.0221438650041818
$ws.Cells.Item(7, 7).Value = 0.0061086523346602
$ws.Cells.Item(7, 8).Value = 0.0325285755097866

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "struggle"
$ws.Cells.Item(8, 3).Value = 0.11651611328125
$ws.Cells.Item(8, 4).Value = -0.4856438636779785
$ws.Cells.Item(8, 5).Value = 0.5658785104751587
$ws.Cells.Item(8, 6).Value = 0.0332921557128429
$ws.Cells.Item(8, 7).Value = -0.0615446716547012
$ws.Cells.Item(8, 8).Value = 0.093156948685646

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "struggle"
$ws.Cells.Item(9, 3).Value = 0.0557413101196289
$ws.Cells.Item(9, 4).Value = 0.3574482798576355
$ws.Cells.Item(9, 5).Value = 0.2321825623512268
$ws.Cells.Item(9, 6).Value = -0.4489859640598297
$ws.Cells.Item(9, 7).Value = -1.353219270706177
$ws.Cells.Item(9, 8).Value = 0.4497495293617248

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "struggle"
$ws.Cells.Item(10, 3).Value = 0.3619680404663086
$ws.Cells.Item(10, 4).Value = 0.0124948024749755
$ws.Cells.Item(10, 5).Value = 0.3587799966335296
$ws.Cells.Item(10, 6).Value = -0.3888157308101654
$ws.Cells.Item(10, 7).Value = -3.63083028793335
$ws.Cells.Item(10, 8).Value = -0.1369865238666534

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "struggle"
$ws.Cells.Item(11, 3).Value = -0.2529764175415039
$ws.Cells.Item(11, 4).Value = 0.1160029470920562
$ws.Cells.Item(11, 5).Value = -0.09882223606109609
$ws.Cells.Item(11, 6).Value = -0.6565274000167847
$ws.Cells.Item(11, 7).Value = -2.371837139129639
$ws.Cells.Item(11, 8).Value = 0.1600466966629028

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "struggle"
$ws.Cells.Item(12, 3).Value = -0.1584005355834961
$ws.Cells.Item(12, 4).Value = 0.0559865832328796
$ws.Cells.Item(12, 5).Value = -0.2031860947608947
$ws.Cells.Item(12, 6).Value = -0.4257730841636657
$ws.Cells.Item(12, 7).Value = -1.438740372657776
$ws.Cells.Item(12, 8).Value = 0.180816113948822

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "struggle"
$ws.Cells.Item(13, 3).Value = -0.1681756973266601
$ws.Cells.Item(13, 4).Value = -0.045459896326065
$ws.Cells.Item(13, 5).Value = 0.3079473972320556
$ws.Cells.Item(13, 6).Value = -0.2063197344541549
$ws.Cells.Item(13, 7).Value = 0.5047274231910706
$ws.Cells.Item(13, 8).Value = -0.1090394482016563

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "struggle"
$ws.Cells.Item(14, 3).Value = 0.7375173568725586
$ws.Cells.Item(14, 4).Value = -0.8549392819404602
$ws.Cells.Item(14, 5).Value = -2.997310400009156
$ws.Cells.Item(14, 6).Value = 1.435685992240906
$ws.Cells.Item(14, 7).Value = 5.099197864532471
$ws.Cells.Item(14, 8).Value = -0.6409503817558289

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "struggle"
$ws.Cells.Item(15, 3).Value = -0.6316938400268555
$ws.Cells.Item(15, 4).Value = 0.0533061251044273
$ws.Cells.Item(15, 5).Value = -1.823783159255981
$ws.Cells.Item(15, 6).Value = 0.7269296646118164
$ws.Cells.Item(15, 7).Value = 4.458247184753418
$ws.Cells.Item(15, 8).Value = 0.2814561724662781

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "struggle"
$ws.Cells.Item(16, 3).Value = 0.1245284080505371
$ws.Cells.Item(16, 4).Value = 0.4134435057640075
$ws.Cells.Item(16, 5).Value = 2.055456638336182
$ws.Cells.Item(16, 6).Value = 0.2370157092809677
$ws.Cells.Item(16, 7).Value = 0.7996225953102112
$ws.Cells.Item(16, 8).Value = 0.1328631937503814

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "struggle"
$ws.Cells.Item(17, 3).Value = -1.905292510986328
$ws.Cells.Item(17, 4).Value = 1.267569422721863
$ws.Cells.Item(17, 5).Value = 0.3008813858032226
$ws.Cells.Item(17, 6).Value = 0.2102903574705124
$ws.Cells.Item(17, 7).Value = 1.452026724815369
$ws.Cells.Item(17, 8).Value = 0.2237294018268585

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "struggle"
$ws.Cells.Item(18, 3).Value = -2.353589773178101
$ws.Cells.Item(18, 4).Value = 0.5766786336898804
$ws.Cells.Item(18, 5).Value = 2.404436111450196
$ws.Cells.Item(18, 6).Value = -0.3715587854385376
$ws.Cells.Item(18, 7).Value = 0.4751004576683044
$ws.Cells.Item(18, 8).Value = 0.1111774742603302

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "struggle"
$ws.Cells.Item(19, 3).Value = -11.09067344665527
$ws.Cells.Item(19, 4).Value = 1.405970811843872
$ws.Cells.Item(19, 5).Value = 10.02403450012207
$ws.Cells.Item(19, 6).Value = 0.1693623960018158
$ws.Cells.Item(19, 7).Value = -1.752572417259216
$ws.Cells.Item(19, 8).Value = 0.1539380401372909

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "struggle"
$ws.Cells.Item(20, 3).Value = 4.286171913146973
$ws.Cells.Item(20, 4).Value = 0.2758489847183227
$ws.Cells.Item(20, 5).Value = -4.509784698486328
$ws.Cells.Item(20, 6).Value = -1.307862520217896
$ws.Cells.Item(20, 7).Value = -5.349499702453613
$ws.Cells.Item(20, 8).Value = -1.575574159622192

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "struggle"
$ws.Cells.Item(21, 3).Value = -1.000519752502441
$ws.Cells.Item(21, 4).Value = -0.010628342628479
$ws.Cells.Item(21, 5).Value = -1.670511245727539
$ws.Cells.Item(21, 6).Value = -0.3645338416099548
$ws.Cells.Item(21, 7).Value = -2.762179851531982
$ws.Cells.Item(21, 8).Value = 0.608421802520752

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "struggle"
$ws.Cells.Item(22, 3).Value = -2.810617446899414
$ws.Cells.Item(22, 4).Value = 0.8466755151748657
$ws.Cells.Item(22, 5).Value = -0.6261429786682129
$ws.Cells.Item(22, 6).Value = -0.3593414723873138
$ws.Cells.Item(22, 7).Value = -2.416888236999512
$ws.Cells.Item(22, 8).Value = -0.4506658315658569

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "struggle"
$ws.Cells.Item(23, 3).Value = -0.552617073059082
$ws.Cells.Item(23, 4).Value = 1.007189750671387
$ws.Cells.Item(23, 5).Value = -2.683732509613037
$ws.Cells.Item(23, 6).Value = -0.3178026378154754
$ws.Cells.Item(23, 7).Value = -1.223715782165527
$ws.Cells.Item(23, 8).Value = -0.2168571650981903

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "struggle"
$ws.Cells.Item(24, 3).Value = -2.832679748535156
$ws.Cells.Item(24, 4).Value = 5.107204437255859
$ws.Cells.Item(24, 5).Value = -6.522222995758057
$ws.Cells.Item(24, 6).Value = 0.042302418500185
$ws.Cells.Item(24, 7).Value = 0.5458080768585205
$ws.Cells.Item(24, 8).Value = 0.195171445608139

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "struggle"
$ws.Cells.Item(25, 3).Value = 0.8469958305358887
$ws.Cells.Item(25, 4).Value = -1.08077871799469
$ws.Cells.Item(25, 5).Value = 7.442714691162109
$ws.Cells.Item(25, 6).Value = 1.255175352096558
$ws.Cells.Item(25, 7).Value = 4.058435916900635
$ws.Cells.Item(25, 8).Value = 0.6265950202941895

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "struggle"
$ws.Cells.Item(26, 3).Value = -3.03963303565979
$ws.Cells.Item(26, 4).Value = 1.802032470703125
$ws.Cells.Item(26, 5).Value = -2.227274417877197
$ws.Cells.Item(26, 6).Value = 1.706299304962158
$ws.Cells.Item(26, 7).Value = 4.895015716552734
$ws.Cells.Item(26, 8).Value = -0.6637051105499268

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "struggle"
$ws.Cells.Item(27, 3).Value = -1.961796522140503
$ws.Cells.Item(27, 4).Value = 1.68219518661499
$ws.Cells.Item(27, 5).Value = 1.394426345825195
$ws.Cells.Item(27, 6).Value = 0.3026837408542633
$ws.Cells.Item(27, 7).Value = 0.5484042763710022
$ws.Cells.Item(27, 8).Value = 0.1058324053883552

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "struggle"
$ws.Cells.Item(28, 3).Value = -2.372189998626709
$ws.Cells.Item(28, 4).Value = 1.225671410560608
$ws.Cells.Item(28, 5).Value = 2.504203796386719
$ws.Cells.Item(28, 6).Value = -0.2704605758190155
$ws.Cells.Item(28, 7).Value = 0.6455318331718445
$ws.Cells.Item(28, 8).Value = -0.1504255682229995

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "struggle"
$ws.Cells.Item(29, 3).Value = -8.473310470581055
$ws.Cells.Item(29, 4).Value = -0.7327957153320312
$ws.Cells.Item(29, 5).Value = 5.200639724731445
$ws.Cells.Item(29, 6).Value = -0.4751004576683044
$ws.Cells.Item(29, 7).Value = -0.1533271819353103
$ws.Cells.Item(29, 8).Value = -0.3729332387447357

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "struggle"
$ws.Cells.Item(30, 3).Value = 3.173869132995605
$ws.Cells.Item(30, 4).Value = -1.535699486732483
$ws.Cells.Item(30, 5).Value = -6.114311695098877
$ws.Cells.Item(30, 6).Value = -0.3394883573055267
$ws.Cells.Item(30, 7).Value = -0.6884451508522034
$ws.Cells.Item(30, 8).Value = -0.1850921660661697

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "struggle"
$ws.Cells.Item(31, 3).Value = -4.002721786499023
$ws.Cells.Item(31, 4).Value = 1.022015571594239
$ws.Cells.Item(31, 5).Value = -0.0432633161544799
$ws.Cells.Item(31, 6).Value = -0.2115120887756347
$ws.Cells.Item(31, 7).Value = -0.2267837226390838
$ws.Cells.Item(31, 8).Value = -0.3090978264808655
